$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header's formatting (bold, border, alignment) onto the
# two new header cells, then overwrite with the new header text.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2..38: I = 1 (constant), J = same value as H (that row)
for ($r = 2; $r -le 38; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
